# The sheet contained five duplicate "0.0" placeholder rows (A6:A10).
# Remove those rows entirely so the remaining data shifts up and the
# now-unused "0.0" shared string is dropped from the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6:A10").EntireRow.Delete()
